# Updates the crypto price/volume table (GitHub Actions style refresh).
# All touched columns (Coin/Link/Price/Volume) are stored as plain text in
# the workbook, so we force a text number format before assigning the new
# value (otherwise Excel would auto-parse strings like "0.490" or
# "1.564.48" as numbers and mangle them), then restore the cell's original
# style so no stray formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextCell 'D2' '28.376.24'
Set-TextCell 'E2' '  -0.41%  '
Set-TextCell 'D3' '1.564.48'
Set-TextCell 'E3' '  -0.20%  '
Set-TextCell 'E4' '  -0.03%  '
Set-TextCell 'D5' '210.87'
Set-TextCell 'E5' '  -0.48%  '
Set-TextCell 'D6' '0.490'
Set-TextCell 'E6' '  -0.58%  '
Set-TextCell 'E7' '  +0.00%  '
Set-TextCell 'D8' '44.53'
Set-TextCell 'E8' '  -3.53%  '
Set-TextCell 'D9' '23.52'
Set-TextCell 'E9' '  -2.20%  '
Set-TextCell 'E10' '  -1.47%  '
Set-TextCell 'D11' '0.0589'
Set-TextCell 'E11' '  -0.66%  '
Set-TextCell 'D12' '0.0895'
Set-TextCell 'E12' '  +1.05%  '
Set-TextCell 'E13' '  -0.20%  '
Set-TextCell 'D14' '1.561.91'
Set-TextCell 'E14' '  -0.40%  '
Set-TextCell 'E15' '  -0.40%  '
Set-TextCell 'D16' '28.355.04'
Set-TextCell 'E16' '  -0.48%  '
Set-TextCell 'E17' '  -1.54%  '
Set-TextCell 'D18' '60.40'
Set-TextCell 'E18' '  -2.95%  '
Set-TextCell 'D19' '228.24'
Set-TextCell 'E19' '  +0.11%  '
Set-TextCell 'D21' '0.0₃0679'
Set-TextCell 'E21' '  -1.92%  '
Set-TextCell 'E22' '  +0.02%  '
Set-TextCell 'E23' '  +1.13%  '
Set-TextCell 'D24' '8.93'
Set-TextCell 'E24' '  -2.11%  '
Set-TextCell 'E25' '  -1.49%  '
Set-TextCell 'D26' '150.20'
Set-TextCell 'E26' '  -0.54%  '
Set-TextCell 'E27' '  -0.88%  '
Set-TextCell 'E28' '  +0.31%  '
Set-TextCell 'E29' '  -2.18%  '
Set-TextCell 'D31' '0.0477'
Set-TextCell 'E31' '  +1.95%  '
Set-TextCell 'E32' '  -3.96%  '
Set-TextCell 'E33' '  -1.14%  '
Set-TextCell 'E34' '  -0.10%  '
Set-TextCell 'D35' '1.386.42'
Set-TextCell 'E35' '  -0.54%  '
Set-TextCell 'E36' '  +1.72%  '
Set-TextCell 'E37' '  -3.52%  '
Set-TextCell 'E38' '  -0.34%  '
Set-TextCell 'D39' '2.65'
Set-TextCell 'E39' '  +2.68%  '
Set-TextCell 'E40' '  -2.07%  '
Set-TextCell 'B41' 'RenderToken'
Set-TextCell 'C41' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D41' '1.95'
Set-TextCell 'E41' '  +3.28%  '
Set-TextCell 'B42' 'ImmutableX'
Set-TextCell 'C42' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D42' '0.518'
Set-TextCell 'E42' '  -3.27%  '
Set-TextCell 'E43' '  +0.00%  '
Set-TextCell 'E44' '  -0.25%  '
Set-TextCell 'E45' '  -1.90%  '
Set-TextCell 'D46' '5.34'
Set-TextCell 'E46' '  -2.97%  '
Set-TextCell 'D47' '0.921'
Set-TextCell 'E47' '  -5.37%  '
Set-TextCell 'D48' '62.17'
Set-TextCell 'E48' '  -1.20%  '
Set-TextCell 'D49' '1.700.50'
Set-TextCell 'E49' '  -0.18%  '
Set-TextCell 'D50' '85.32'
Set-TextCell 'E50' '  -0.68%  '
Set-TextCell 'B51' 'BabyDogeCoin'
Set-TextCell 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D51' '0.0₆0100'
Set-TextCell 'E51' '  -1.97%  '
